# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" tracking sheet and
# moves the "last row" date formatting down to the newly appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (29) used a distinct "date only" number format to mark
# it as the latest entry. Since a new row is being appended below it, that
# row reverts to the regular "date + time" number format used by all other
# data rows.
$ws.Range("A29").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 30.
$ws.Range("A30").Value = 45614
$ws.Range("B30").Value = 76
$ws.Range("C30").Value = 62
$ws.Range("D30").Value = 73

# The newly appended row becomes the new "last row" and gets the distinct
# date-only number format previously on row 29.
$ws.Range("A30").NumberFormat = "YYYY-MM-DD"
